$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows with new data
# Price column is forced to Text format to preserve the original
# non-numeric "1.234.56"-style formatting and trailing zeros, exactly
# like the source inline strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.972.55"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.747.90"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9981"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.04"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9982"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5170"
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2816"
$ws.Range("E8").Value = "  +7.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.72"
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06136"
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.742.70"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07017"
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.47"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6436"
$ws.Range("E14").Value = "  +6.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.526"
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.10"
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9976"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9979"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "25.965.64"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.51"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006630"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.963.02"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.145"
$ws.Range("E23").Value = "  +2.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.579"
$ws.Range("E24").Value = "  +4.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.158"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.36"
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.503"
$ws.Range("E27").Value = "  +2.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.843"
$ws.Range("E28").Value = "  +1.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "103.00"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08318"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.644"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.435"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04419"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.615"
$ws.Range("E35").Value = "  -1.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9861"
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6116"
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.678"
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01573"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.940"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9975"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.75"
$ws.Range("E42").Value = "  -2.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3875"
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7352"
$ws.Range("E44").Value = "  -3.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.008"
$ws.Range("E45").Value = "  +3.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05460"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.385"
$ws.Range("E47").Value = "  +7.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1121"
$ws.Range("E48").Value = "  +3.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.76"
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "29.91"
$ws.Range("E50").Value = "  -0.74%  "

# Row 51: coin entry was swapped from EnergySwap to Decentraland
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3428"
$ws.Range("E51").Value = "  +0.08%  "
